$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 26/27: coin identity (name + link) swapped, plus new price/volume ---
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'21.86"
$ws.Range("E26").Value = "  +4.60%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.115.53"
$ws.Range("E27").Value = "  +0.05%  "

# --- Price / Volume(1h) updates for remaining rows ---
$ws.Range("D2").Value = "30.159.44"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "1.897.38"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'325.71"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.5180"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "'0.4010"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.08443"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'42.78"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "'1.119"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "'23.32"
$ws.Range("E12").Value = "  +12.91%  "
$ws.Range("D13").Value = "'6.451"
$ws.Range("E13").Value = "  +3.33%  "
$ws.Range("D14").Value = "1.893.28"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "'94.77"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").Value = "'0.06661"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "'18.28"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'5.949"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").Value = "30.177.27"
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("D24").Value = "'11.27"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D28").Value = "'161.19"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Value = "'2.389"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "'129.20"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "'1.093"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").Value = "'0.1058"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'6.056"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").Value = "'3.703"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").Value = "'0.02495"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "'0.06562"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "'5.253"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").Value = "'1.215"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").Value = "'11.80"
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").Value = "'8.777"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'1.234"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "'0.6110"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "'3.709"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D48").Value = "'1.240"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "'124.57"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").Value = "'1.163"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  +1.91%  "
